$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 values (stat_sample_size changed from 100 to 1000,
#     and the area estimates were recomputed accordingly) ---
$ws.Range("B2").Value = 1.50922833152
$ws.Range("C2").Value = 1.50922833152
$ws.Range("D2").Value = 1.5106093184
$ws.Range("E2").Value = 1.5103963056
$ws.Range("F2").Value = 1.510551224
$ws.Range("G2").Value = 1000

# --- Update existing row 3 values (var_area row) ---
$ws.Range("B3").Value = 0.0006231141841015626
$ws.Range("C3").Value = 0.0006231141841015626
$ws.Range("D3").Value = 0.00001297575825395703
$ws.Range("E3").Value = 0.0002493197787103174
$ws.Range("F3").Value = 0.0009127010069212253
$ws.Range("G3").Value = 1000

# --- Add new row 4: ci_upper ---
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = "ci_upper"
$ws.Range("B4").Value = 1.558212799840014
$ws.Range("C4").Value = 1.558212799840014
$ws.Range("D4").Value = 1.517678041226865
$ws.Range("E4").Value = 1.541381408750631
$ws.Range("F4").Value = 1.569835407754516
$ws.Range("G4").Value = 1000

# --- Add new row 5: ci_down ---
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "ci_down"
$ws.Range("B5").Value = 1.460243863199986
$ws.Range("C5").Value = 1.460243863199986
$ws.Range("D5").Value = 1.503540595573135
$ws.Range("E5").Value = 1.47941120244937
$ws.Range("F5").Value = 1.451267040245484
$ws.Range("G5").Value = 1000
